# Weekly update: a new price observation (week of 2022-06-17) is inserted
# at row 21, pushing the existing rows 21-39 down to 22-40.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 21 (shifts rows 21:39 -> 22:40, carrying
# their formatting/styles down with them, and extends the sheet dimension).
$ws.Rows("21:21").Insert()

# Populate the newly inserted row 21 with the new weekly record.
$ws.Range("A21").Value = 9
$ws.Range("B21").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C21").Value = "Metropolitana"
$ws.Range("D21").Value = 44729
$ws.Range("E21").Value = 13
$ws.Range("F21").Value = 100112035
$ws.Range("G21").Value = "Bruselas (repollito)"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 52
$ws.Range("K21").Value = 24000
$ws.Range("L21").Value = 24000
$ws.Range("M21").Value = 24000
$ws.Range("N21").Value = "$/malla 15 kilos"
$ws.Range("O21").Value = "Hijuelas"
$ws.Range("P21").Value = 1600
$ws.Range("Q21").Value = 15
$ws.Range("R21").Value = "Hortaliza"
